$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.227.06"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.689.25"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.29%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0892"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.926.10"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "1.677.40"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "27.214.80"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("E23").Value = "  +5.19%  "
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "1.583.95"
$ws.Range("E32").Value = "  +6.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.604"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "1.835.56"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +5.85%  "
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.72%  "
$ws.Range("E51").Value = "  +3.18%  "
